# Update PlanID (column F) values for the set of countries whose
# plan mapping changed in this data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(35, 6).Value = 22
$ws.Cells.Item(51, 6).Value = 26
$ws.Cells.Item(54, 6).Value = 26
$ws.Cells.Item(58, 6).Value = 26
$ws.Cells.Item(61, 6).Value = 26
$ws.Cells.Item(64, 6).Value = 26
$ws.Cells.Item(66, 6).Value = 26
$ws.Cells.Item(68, 6).Value = 26
$ws.Cells.Item(70, 6).Value = 26
$ws.Cells.Item(72, 6).Value = 26
$ws.Cells.Item(74, 6).Value = 26
$ws.Cells.Item(76, 6).Value = 26
$ws.Cells.Item(77, 6).Value = 26
$ws.Cells.Item(78, 6).Value = 26
$ws.Cells.Item(82, 6).Value = 26
$ws.Cells.Item(83, 6).Value = 26
$ws.Cells.Item(84, 6).Value = 26
$ws.Cells.Item(86, 6).Value = 26
$ws.Cells.Item(87, 6).Value = 26
$ws.Cells.Item(88, 6).Value = 26
$ws.Cells.Item(90, 6).Value = 26
$ws.Cells.Item(92, 6).Value = 26
$ws.Cells.Item(110, 6).Value = 14
$ws.Cells.Item(128, 6).Value = 26
$ws.Cells.Item(133, 6).Value = 30
$ws.Cells.Item(135, 6).Value = 30
$ws.Cells.Item(138, 6).Value = 30
$ws.Cells.Item(139, 6).Value = 30
$ws.Cells.Item(140, 6).Value = 30
$ws.Cells.Item(141, 6).Value = 30
$ws.Cells.Item(144, 6).Value = 30
$ws.Cells.Item(147, 6).Value = 30
$ws.Cells.Item(150, 6).Value = 29
$ws.Cells.Item(152, 6).Value = 30
$ws.Cells.Item(157, 6).Value = 29
$ws.Cells.Item(161, 6).Value = 30
$ws.Cells.Item(162, 6).Value = 29
$ws.Cells.Item(163, 6).Value = 30
$ws.Cells.Item(164, 6).Value = 30
$ws.Cells.Item(168, 6).Value = 29
$ws.Cells.Item(171, 6).Value = 29
$ws.Cells.Item(172, 6).Value = 30
$ws.Cells.Item(173, 6).Value = 29
$ws.Cells.Item(234, 6).Value = 26
$ws.Cells.Item(235, 6).Value = 30
$ws.Cells.Item(240, 6).Value = 36
$ws.Cells.Item(241, 6).Value = 26
$ws.Cells.Item(242, 6).Value = 26
$ws.Cells.Item(243, 6).Value = 26
$ws.Cells.Item(244, 6).Value = 26
$ws.Cells.Item(245, 6).Value = 26
$ws.Cells.Item(246, 6).Value = 26
$ws.Cells.Item(247, 6).Value = 26
$ws.Cells.Item(248, 6).Value = 26
$ws.Cells.Item(249, 6).Value = 26
$ws.Cells.Item(250, 6).Value = 26
$ws.Cells.Item(251, 6).Value = 26
$ws.Cells.Item(252, 6).Value = 26
$ws.Cells.Item(253, 6).Value = 26
$ws.Cells.Item(254, 6).Value = 26
$ws.Cells.Item(255, 6).Value = 26
$ws.Cells.Item(256, 6).Value = 26
$ws.Cells.Item(257, 6).Value = 26
$ws.Cells.Item(258, 6).Value = 26
$ws.Cells.Item(259, 6).Value = 26
$ws.Cells.Item(260, 6).Value = 26
$ws.Cells.Item(261, 6).Value = 26
$ws.Cells.Item(262, 6).Value = 26
$ws.Cells.Item(263, 6).Value = 26
$ws.Cells.Item(264, 6).Value = 37
$ws.Cells.Item(265, 6).Value = 10
$ws.Cells.Item(266, 6).Value = 16
$ws.Cells.Item(267, 6).Value = 23
$ws.Cells.Item(268, 6).Value = 26
$ws.Cells.Item(269, 6).Value = 26
$ws.Cells.Item(270, 6).Value = 26
$ws.Cells.Item(271, 6).Value = 26
$ws.Cells.Item(272, 6).Value = 26
$ws.Cells.Item(273, 6).Value = 35
$ws.Cells.Item(278, 6).Value = 12
$ws.Cells.Item(281, 6).Value = 33
$ws.Cells.Item(282, 6).Value = 24
$ws.Cells.Item(293, 6).Value = 27
$ws.Cells.Item(294, 6).Value = 32

# Column H (PhoneCode) now has real content worth auto-sizing - match the
# author's "best fit" column-width action.
$ws.Columns.Item(8).EntireColumn.AutoFit() | Out-Null

# Restore the author's last-used cell selection on the sheet.
$ws.Range("J17").Select() | Out-Null
